# Update the keyboard-key-spacing measurements (new "aktuelle Koordinaten
# der neuen Tastatur") and append a "Leertaste" (spacebar) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Re-measured numeric values (mm offsets / widths) for rows 3..8.
#    The letter labels in between (B3, D3, F3, ... etc.) are untouched.
# ---------------------------------------------------------------------
$values = [ordered]@{
    "C3" = -135
    "E3" = -106
    "G3" = -77
    "I3" = -48
    "K3" = -19
    "M3" = 10
    "O3" = 39
    "Q3" = 68
    "S3" = 97
    "U3" = 126
    "V3" = 30
    "B4" = 213
    "D4" = 213
    "F4" = 213
    "H4" = 213
    "J4" = 213
    "L4" = 213
    "N4" = 213
    "P4" = 213
    "R4" = 213
    "T4" = 213
    "C5" = -106
    "E5" = -77
    "G5" = -48
    "I5" = -19
    "K5" = 10
    "M5" = 39
    "O5" = 68
    "Q5" = 97
    "S5" = 126
    "U5" = 155
    "V5" = 30
    "B6" = 188
    "D6" = 188
    "F6" = 188
    "H6" = 188
    "J6" = 188
    "L6" = 188
    "N6" = 188
    "P6" = 188
    "R6" = 188
    "T6" = 188
    "C7" = -106
    "E7" = -77
    "G7" = -48
    "I7" = -19
    "K7" = 10
    "M7" = 39
    "O7" = 68
    "V7" = 30
    "B8" = 163
    "D8" = 163
    "F8" = 163
    "H8" = 163
    "J8" = 163
    "L8" = 163
    "N8" = 163
}

foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# ---------------------------------------------------------------------
# 2) The measurement cells in rows 3/5/7 (fill colour "2") and rows
#    4/6/8 (fill colour "3") used to carry an extra left-alignment on
#    top of their fill. That alignment was removed - re-apply the
#    already-used plain-fill format (copied from a sibling cell that
#    already has the target look) so the workbook doesn't grow new,
#    duplicate cell styles.
# ---------------------------------------------------------------------
$plainFill2 = $ws.Range("C2")   # fill "2", no alignment override
$plainFill3 = $ws.Range("A4")   # fill "3", no alignment override

$fill2Cells = @("C3","E3","G3","I3","K3","M3","O3","Q3","S3","U3", `
                "C5","E5","G5","I5","K5","M5","O5","Q5","S5","U5", `
                "C7","E7","G7","I7","K7","M7","O7")
$fill3Cells = @("B4","D4","F4","H4","J4","L4","N4","P4","R4","T4", `
                "B6","D6","F6","H6","J6","L6","N6","P6","R6","T6", `
                "B8","D8","F8","H8","J8","L8","N8")

$plainFill2.Copy()
foreach ($ref in $fill2Cells) {
    $ws.Range($ref).PasteSpecial(-4122)
}

$plainFill3.Copy()
foreach ($ref in $fill3Cells) {
    $ws.Range($ref).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) New "Leertaste" (spacebar) label below the keyboard layout.
# ---------------------------------------------------------------------
$ws.Range("B10").Value = "Leertaste"

# ---------------------------------------------------------------------
# 4) Restore the active selection to where editing left off.
# ---------------------------------------------------------------------
$ws.Range("N8").Select()
